$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores plain text (e.g. "30.574.65", "0.9999") rather than
# numbers. For updated prices that look numeric, force the cell format to Text
# first so Excel does not silently convert them into real numbers.
$priceCellsNeedingTextFormat = @(
    "D4", "D5", "D7", "D8", "D9", "D10", "D11", "D13", "D14", "D15",
    "D16", "D18", "D19", "D20", "D21", "D23", "D24", "D25", "D26", "D27",
    "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D39",
    "D41", "D42", "D43", "D45", "D47", "D48", "D50", "D51"
)
foreach ($cellRef in $priceCellsNeedingTextFormat) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Refresh each coin row with the latest scraped price / 1h-volume figures.

$ws.Range("D2").Value = "30.574.65"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.934.52"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "246.00"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "0.4861"
$ws.Range("E7").Value = "  +2.73%  "
$ws.Range("D8").Value = "0.2917"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "0.06809"
$ws.Range("E9").Value = "  -0.70%  "
$ws.Range("D10").Value = "112.72"
$ws.Range("E10").Value = "  +6.30%  "
$ws.Range("D11").Value = "19.49"
$ws.Range("E11").Value = "  +5.49%  "
$ws.Range("D12").Value = "1.936.80"
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").Value = "5.515"
$ws.Range("E13").Value = "  +3.19%  "
$ws.Range("D14").Value = "0.07591"
$ws.Range("E14").Value = "  -1.81%  "
$ws.Range("D15").Value = "0.6807"
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("D16").Value = "299.32"
$ws.Range("E16").Value = "  +3.44%  "
$ws.Range("D17").Value = "30.570.34"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").Value = "13.13"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("D19").Value = "0.000007669"
$ws.Range("E19").Value = "  +0.18%  "
# Dai and Uniswap swapped rank positions in this run.
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "0.9999"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "5.562"
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").Value = "2.186.43"
$ws.Range("E22").Value = "  +0.29%  "
$ws.Range("D23").Value = "0.9998"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "6.520"
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("D25").Value = "9.572"
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("D26").Value = "168.24"
$ws.Range("E26").Value = "  +0.83%  "
$ws.Range("D27").Value = "20.43"
$ws.Range("E27").Value = "  -1.87%  "
$ws.Range("D28").Value = "2.130"
$ws.Range("E28").Value = "  -0.22%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").Value = "1.434"
$ws.Range("E30").Value = "  +1.95%  "
$ws.Range("D31").Value = "4.183"
$ws.Range("E31").Value = "  -0.60%  "
$ws.Range("D32").Value = "4.103"
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("D33").Value = "0.05011"
$ws.Range("E33").Value = "  -0.85%  "
$ws.Range("D34").Value = "0.7497"
$ws.Range("E34").Value = "  +1.95%  "
$ws.Range("D35").Value = "1.148"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  -1.02%  "
$ws.Range("D37").Value = "2.716"
$ws.Range("E37").Value = "  -1.06%  "
$ws.Range("D38").Value = "2.693"
$ws.Range("E38").Value = "  +0.36%  "
$ws.Range("D39").Value = "2.034"
$ws.Range("E39").Value = "  -1.14%  "
$ws.Range("E40").Value = "  -1.33%  "
$ws.Range("D41").Value = "0.4465"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").Value = "0.8733"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").Value = "5.838"
$ws.Range("E43").Value = "  -1.35%  "
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").Value = "69.79"
$ws.Range("E45").Value = "  +2.53%  "
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("D47").Value = "49.23"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").Value = "9.319"
$ws.Range("E48").Value = "  -0.88%  "
$ws.Range("E49").Value = "  -2.30%  "
$ws.Range("D50").Value = "0.2546"
$ws.Range("E50").Value = "  +2.83%  "
$ws.Range("D51").Value = "35.11"
$ws.Range("E51").Value = "  -0.54%  "
